$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -eq "⬛") {
            $cell.Value = "📘"
        }
        elseif ($v -eq "noir") {
            $cell.Value = "bleu"
        }
    }
}
